$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 929.6667
$ws.Range("I19").Value = 915.5
$ws.Range("K19").Value = 915.5
$ws.Range("M19").Value = -740.5
$ws.Range("H21").Value = 1000
$ws.Range("I21").Value = 1000
$ws.Range("K21").Value = 1000
$ws.Range("M21").Value = -532
$ws.Range("H23").Value = 1000
$ws.Range("I23").Value = 1000
$ws.Range("K23").Value = 1000
$ws.Range("M23").Value = -766
$ws.Range("H33").Value = 494.25
$ws.Range("I33").Value = 590
$ws.Range("K33").Value = 590
$ws.Range("M33").Value = -361
$ws.Range("H39").Value = 195.3
$ws.Range("I39").Value = 152.625
$ws.Range("J39").Value = 366
$ws.Range("K39").Value = 457.875
$ws.Range("L39").Value = 1098
$ws.Range("M39").Value = -161.875
$ws.Range("N39").Value = -1690
$ws.Range("H80").Value = 3068
$ws.Range("I80").Value = 2122.8
$ws.Range("J80").Value = 3927.2727
$ws.Range("K80").Value = 6368.400000000001
$ws.Range("L80").Value = 11781.8181
$ws.Range("M80").Value = -5370.400000000001
$ws.Range("N80").Value = -13777.8181
$ws.Range("H83").Value = 3068
$ws.Range("I83").Value = 2122.8
$ws.Range("J83").Value = 3927.2727
$ws.Range("K83").Value = 19105.2
$ws.Range("L83").Value = 35345.4543
$ws.Range("M83").Value = -14113.2
$ws.Range("N83").Value = -45329.4543
$ws.Range("H96").Value = 1477.5
$ws.Range("I96").Value = 1628.75
$ws.Range("K96").Value = 4886.25
$ws.Range("M96").Value = -3513.25
$ws.Range("H101").Value = 1327.2727
$ws.Range("J101").Value = 1877.6
$ws.Range("L101").Value = 5632.799999999999
$ws.Range("N101").Value = -8876.799999999999
$ws.Range("H119").Value = 999.5
$ws.Range("J119").Value = 999.5
$ws.Range("L119").Value = 2998.5
$ws.Range("N119").Value = -12674.5
$ws.Range("H135").Value = 2313.6428
$ws.Range("I135").Value = 1032.5
$ws.Range("J135").Value = 10000.5
$ws.Range("K135").Value = 9292.5
$ws.Range("L135").Value = 90004.5
$ws.Range("M135").Value = -6757.5
$ws.Range("N135").Value = -95074.5
$ws.Range("H138").Value = 3624.1562
$ws.Range("I138").Value = 7399
$ws.Range("J138").Value = 3233.6553
$ws.Range("K138").Value = 22197
$ws.Range("L138").Value = 9700.965899999999
$ws.Range("M138").Value = -17057
$ws.Range("N138").Value = -19980.9659

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 272.4
$ws.Range("I5").Value = 109
$ws.Range("J5").Value = 381.33334
$ws.Range("K5").Value = 109
$ws.Range("L5").Value = 381.33334
$ws.Range("M5").Value = 3
$ws.Range("N5").Value = -605.33334
$ws.Range("H32").Value = 1830.9592
$ws.Range("I32").Value = 1282.4
$ws.Range("J32").Value = 8002.25
$ws.Range("K32").Value = 1282.4
$ws.Range("L32").Value = 8002.25
$ws.Range("M32").Value = -995.4000000000001
$ws.Range("N32").Value = -8576.25
$ws.Range("H102").Value = 1665.25
$ws.Range("I102").Value = 1692.5555
$ws.Range("K102").Value = 1692.5555
$ws.Range("M102").Value = -70.55549999999994

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 272.4
$ws.Range("I4").Value = 109
$ws.Range("J4").Value = 381.33334
$ws.Range("K4").Value = 109
$ws.Range("L4").Value = 381.33334
$ws.Range("M4").Value = 6
$ws.Range("N4").Value = -611.33334
$ws.Range("H19").Value = 8000
$ws.Range("I19").Value = 8000
$ws.Range("K19").Value = 8000
$ws.Range("M19").Value = -7827
$ws.Range("H105").Value = 28648
$ws.Range("I105").Value = 28909.625
$ws.Range("K105").Value = 28909.625
$ws.Range("M105").Value = -27162.625
$ws.Range("H134").Value = 1654.8572
$ws.Range("I134").Value = 1160.0741
$ws.Range("J134").Value = 15014
$ws.Range("K134").Value = 3480.2223
$ws.Range("L134").Value = 45042
$ws.Range("M134").Value = -945.2223000000004
$ws.Range("N134").Value = -50112

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2735.6858
$ws.Range("I58").Value = 1210.8334
$ws.Range("J58").Value = 11884.8
$ws.Range("K58").Value = 1210.8334
$ws.Range("L58").Value = 11884.8
$ws.Range("M58").Value = -1007.8334
$ws.Range("N58").Value = -12290.8
$ws.Range("H86").Value = 16001
$ws.Range("I86").Value = 12500
$ws.Range("K86").Value = 12500
$ws.Range("M86").Value = -11377
$ws.Range("H89").Value = 16001
$ws.Range("I89").Value = 12500
$ws.Range("K89").Value = 62500
$ws.Range("M89").Value = -56884
$ws.Range("H105").Value = 7747.5454
$ws.Range("J105").Value = 5877.5
$ws.Range("L105").Value = 5877.5
$ws.Range("N105").Value = -9371.5
$ws.Range("H132").Value = 2012.3462
$ws.Range("I132").Value = 1332.28
$ws.Range("J132").Value = 19014
$ws.Range("K132").Value = 3996.84
$ws.Range("L132").Value = 57042
$ws.Range("M132").Value = -1466.84
$ws.Range("N132").Value = -62102
$ws.Range("H136").Value = 2735.6858
$ws.Range("I136").Value = 1210.8334
$ws.Range("J136").Value = 11884.8
$ws.Range("K136").Value = 3632.5002
$ws.Range("L136").Value = 35654.39999999999
$ws.Range("M136").Value = -1082.5002
$ws.Range("N136").Value = -40754.39999999999

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H42").Value = 13001.333
$ws.Range("J42").Value = 13001.333
$ws.Range("L42").Value = 39003.999
$ws.Range("N42").Value = -40071.999
$ws.Range("H69").Value = 8008.3335
$ws.Range("I69").Value = 5012
$ws.Range("K69").Value = 15036
$ws.Range("M69").Value = -14225
$ws.Range("H72").Value = 8008.3335
$ws.Range("I72").Value = 5012
$ws.Range("K72").Value = 45108
$ws.Range("M72").Value = -41052
$ws.Range("H121").Value = 3434.75
$ws.Range("J121").Value = 3795.6
$ws.Range("L121").Value = 11386.8
$ws.Range("N121").Value = -14006.8
$ws.Range("H122").Value = 1441.7742
$ws.Range("I122").Value = 560.8182
$ws.Range("J122").Value = 1926.3
$ws.Range("K122").Value = 5047.3638
$ws.Range("L122").Value = 17336.7
$ws.Range("M122").Value = -2597.3638
$ws.Range("N122").Value = -22236.7

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 7704.909
$ws.Range("I80").Value = 8050
$ws.Range("J80").Value = 7417.3335
$ws.Range("K80").Value = 8050
$ws.Range("L80").Value = 7417.3335
$ws.Range("M80").Value = -7052
$ws.Range("N80").Value = -9413.333500000001
$ws.Range("H83").Value = 7704.909
$ws.Range("I83").Value = 8050
$ws.Range("J83").Value = 7417.3335
$ws.Range("K83").Value = 40250
$ws.Range("L83").Value = 37086.6675
$ws.Range("M83").Value = -35258
$ws.Range("N83").Value = -47070.6675

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 5099.5
$ws.Range("I46").Value = 795
$ws.Range("J46").Value = 5714.4287
$ws.Range("K46").Value = 795
$ws.Range("L46").Value = 5714.4287
$ws.Range("M46").Value = -607
$ws.Range("N46").Value = -6090.4287
$ws.Range("H68").Value = 4232.9165
$ws.Range("I68").Value = 2979.8
$ws.Range("J68").Value = 5128
$ws.Range("K68").Value = 2979.8
$ws.Range("L68").Value = 5128
$ws.Range("M68").Value = -2230.8
$ws.Range("N68").Value = -6626
$ws.Range("H71").Value = 4232.9165
$ws.Range("I71").Value = 2979.8
$ws.Range("J71").Value = 5128
$ws.Range("K71").Value = 14899
$ws.Range("L71").Value = 25640
$ws.Range("M71").Value = -11155
$ws.Range("N71").Value = -33128

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H82").Value = 50000
$ws.Range("J82").Value = 0
$ws.Range("L82").Value = 0
$ws.Range("H85").Value = 50000
$ws.Range("J85").Value = 0
$ws.Range("L85").Value = 0
$ws.Range("N82").ClearContents()
$ws.Range("N85").ClearContents()

Write-Host "Applied all Lamia_Profits updates"